# Generate Report for Handback
#
# The "46bd691b-2c68-4923-a222-5e3a92f78283.md" file now fails the handback
# transform, so its status flips from "Ready for handoff" to
# "Handback transform failed" everywhere it is reported (Overview + the
# per-locale detail sheets), and the per-locale sheets get a detailed error
# message in their "Error Detail" column (with that column widened so the
# message is readable).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row 3 is the 46bd691b...md file, zh-cn (E) and de-de (F) columns.
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# zh-cn detail sheet: row 3 is the same file; Status column C, Error Detail column P.
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("P3").Value = "Handback file name: ouhdlfq0.4br is different with handoff file name: 46bd691b-2c68-4923-a222-5e3a92f78283.1989e98e9b6dd71f79697b086c2222feea3b0904.zh-cn."

# de-de detail sheet: row 3 is the same file; Status column C, Error Detail column P.
$dede.Range("C3").Value = $newStatus
$dede.Range("P3").Value = "Handback file name: ouhdlfq0.4br is different with handoff file name: 46bd691b-2c68-4923-a222-5e3a92f78283.1989e98e9b6dd71f79697b086c2222feea3b0904.de-de."

# Widen the "Error Detail" column (P) on both detail sheets to fit the new
# message, matching the width already used by column A (40 chars).
$refWidth = $zhcn.Range("A1").ColumnWidth()
$zhcn.Range("P1:P3").ColumnWidth = $refWidth
$dede.Range("P1:P3").ColumnWidth = $refWidth
